$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.454.86'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '3.503.83'
$ws.Range("E3").Value = '  -0.26%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '607.08'
$ws.Range("E5").Value = '  +4.73%  '
$ws.Range("D6").Value = '169.59'
$ws.Range("E6").Value = '  -2.91%  '
$ws.Range("D7").Value = '0.617'
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("D8").Value = '3.504.57'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("D10").Value = '0.196'
$ws.Range("E10").Value = '  +3.89%  '
$ws.Range("E11").Value = '  -0.64%  '
$ws.Range("D12").Value = '0.578'
$ws.Range("E12").Value = '  -3.75%  '
$ws.Range("D13").Value = '46.95'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000278'
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").Value = '4.071.30'
$ws.Range("E15").Value = '  -0.09%  '
$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '612.80'
$ws.Range("E16").Value = '  -8.81%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.30'
$ws.Range("E17").Value = '  -6.16%  '
$ws.Range("D18").Value = '3.502.20'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").Value = '69.461.29'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").Value = '17.21'
$ws.Range("E21").Value = '  -2.08%  '
$ws.Range("D22").Value = '10.15'
$ws.Range("E22").Value = '  -9.50%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.880'
$ws.Range("E23").Value = '  -2.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.80'
$ws.Range("E24").Value = '  -3.07%  '
$ws.Range("D25").Value = '95.59'
$ws.Range("E25").Value = '  -2.83%  '
$ws.Range("D26").Value = '3.86'
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '2.59'
$ws.Range("E28").Value = '  -2.34%  '
$ws.Range("D29").Value = '9.19'
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("D30").Value = '33.12'
$ws.Range("E30").Value = '  +0.43%  '
$ws.Range("D31").Value = '8.41'
$ws.Range("E31").Value = '  -3.96%  '
$ws.Range("E32").Value = '  -4.57%  '
$ws.Range("D33").Value = '1.33'
$ws.Range("E33").Value = '  -2.41%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.90'
$ws.Range("E34").Value = '  -5.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '555.70'
$ws.Range("E35").Value = '  -3.83%  '
$ws.Range("D36").Value = '10.74'
$ws.Range("E36").Value = '  -1.95%  '
$ws.Range("D37").Value = '3.49'
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("D38").Value = '56.75'
$ws.Range("E38").Value = '  -0.50%  '
$ws.Range("D39").Value = '0.101'
$ws.Range("E39").Value = '  -4.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").Value = '0.0448'
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("D42").Value = '0.139'
$ws.Range("E42").Value = '  +1.51%  '
$ws.Range("D43").Value = '3.332.73'
$ws.Range("E43").Value = '  -2.65%  '
$ws.Range("D44").Value = '0.325'
$ws.Range("E44").Value = '  -3.76%  '
$ws.Range("D45").Value = '32.96'
$ws.Range("E45").Value = '  -1.46%  '
$ws.Range("D46").Value = '0.0₃0697'
$ws.Range("E46").Value = '  -1.37%  '
$ws.Range("D47").Value = '2.89'
$ws.Range("E47").Value = '  -0.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("D49").Value = '0.129'
$ws.Range("E49").Value = '  -3.38%  '
$ws.Range("D50").Value = '135.65'
$ws.Range("E50").Value = '  +3.71%  '
$ws.Range("E51").Value = '  +7.37%  '
